$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Cat"
$ws.Range("B3").Value = "Cat"
$ws.Range("C3").Value = "Katt"
$ws.Range("D3").Value = "Need review"

$ws.Range("A4").Value = "Password"
$ws.Range("B4").Value = "Password"
$ws.Range("C4").Value = "Passord"
$ws.Range("D4").Value = "Need review"

$ws.Range("A5").Value = "Run if you like"
$ws.Range("B5").Value = "Run if you like"
$ws.Range("C5").Value = "Løp hvis du vil"
$ws.Range("D5").Value = "Need review"

$ws.Range("A6").Value = "House"
$ws.Range("B6").Value = "House"
$ws.Range("C6").Value = "Hus"
$ws.Range("D6").Value = "Need review"

$ws.Range("A7").Value = "Wood"
$ws.Range("B7").Value = "Wood"
$ws.Range("C7").Value = "Tre"
$ws.Range("D7").Value = "Need review"
